# Updated symbol list on Mon Jan  9 21:26:10 UTC 2023 with GitHub Actions
#
# All price/volume cells in this sheet are stored as plain text (not
# numbers/percentages), so every new value is written with a leading
# apostrophe to force Excel to keep it as literal text instead of
# auto-converting it to a Number/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'1.59%"
$ws.Range("D3").Value = "'26.73"
$ws.Range("E3").Value = "'0.09%"
$ws.Range("D4").Value = "'4.908"
$ws.Range("E4").Value = "'4.15%"
$ws.Range("D5").Value = "'0.06330"
$ws.Range("E5").Value = "'3.79%"
$ws.Range("D6").Value = "'6.920"
$ws.Range("E6").Value = "'2.57%"
$ws.Range("D7").Value = "'3.356"
$ws.Range("E7").Value = "'5.77%"
$ws.Range("D8").Value = "'1.348"
$ws.Range("E8").Value = "'51.59%"
$ws.Range("D9").Value = "'0.8841"
$ws.Range("E9").Value = "'3.30%"
$ws.Range("D10").Value = "'0.1474"
$ws.Range("E10").Value = "'3.06%"
$ws.Range("D11").Value = "'0.05072"
$ws.Range("E11").Value = "'2.47%"
$ws.Range("D12").Value = "'0.07401"
$ws.Range("E12").Value = "'4.14%"
$ws.Range("D13").Value = "'0.03182"
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("D14").Value = "'0.09049"
$ws.Range("E14").Value = "'0.13%"
$ws.Range("D15").Value = "'0.001550"
$ws.Range("E15").Value = "'0.77%"
$ws.Range("D16").Value = "'0.0006292"
$ws.Range("E16").Value = "'3.43%"
$ws.Range("D17").Value = "'0.006071"
$ws.Range("E17").Value = "'1.97%"
$ws.Range("D18").Value = "'3.472"
$ws.Range("E18").Value = "'0.26%"
$ws.Range("D19").Value = "'2.283"
$ws.Range("E19").Value = "'0.83%"
$ws.Range("E21").Value = "'2.77%"
$ws.Range("D22").Value = "'3.926"
$ws.Range("E22").Value = "'2.32%"
$ws.Range("D23").Value = "'0.04329"
$ws.Range("E23").Value = "'1.83%"
$ws.Range("D24").Value = "'0.001174"
$ws.Range("E24").Value = "'-0.25%"
$ws.Range("D25").Value = "'0.003647"
$ws.Range("E25").Value = "'-12.09%"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("D27").Value = "'0.0001701"
$ws.Range("E27").Value = "'1.16%"
$ws.Range("D40").Value = "'0.04046"
$ws.Range("E40").Value = "'2.43%"
$ws.Range("D41").Value = "'0.006636"
$ws.Range("E41").Value = "'58.28%"
$ws.Range("D42").Value = "'0.1168"
$ws.Range("E42").Value = "'4.26%"
$ws.Range("D43").Value = "'0.002238"
$ws.Range("E43").Value = "'9.88%"
$ws.Range("D44").Value = "'0.01263"
$ws.Range("E44").Value = "'7.11%"
$ws.Range("D45").Value = "'0.00005334"
$ws.Range("E45").Value = "'3.86%"
$ws.Range("E46").Value = "'141.64%"
$ws.Range("E47").Value = "'-13.22%"
